# Scheduled data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H,I,J,K,L,M,N) across all Leve sheets with latest market values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 635.3333
$ws.Range("I33").Value = 340
$ws.Range("K33").Value = 340
$ws.Range("M33").Value = -111
$ws.Range("H69").Value = 11714.286
$ws.Range("J69").Value = 12333.333
$ws.Range("L69").Value = 36999.999
$ws.Range("N69").Value = -38747.999
$ws.Range("H72").Value = 11714.286
$ws.Range("J72").Value = 12333.333
$ws.Range("L72").Value = 110999.997
$ws.Range("N72").Value = -119735.997
$ws.Range("H112").Value = 7005.5527
$ws.Range("J112").Value = 7005.5527
$ws.Range("L112").Value = 21016.6581
$ws.Range("N112").Value = -23232.6581
$ws.Range("H133").Value = 89861.55499999999
$ws.Range("J133").Value = 89861.55499999999
$ws.Range("L133").Value = 89861.55499999999
$ws.Range("N133").Value = -99981.55499999999
$ws.Range("H134").Value = 89938.11
$ws.Range("J134").Value = 89938.11
$ws.Range("L134").Value = 89938.11
$ws.Range("N134").Value = -100078.11
$ws.Range("H135").Value = 2114.7917
$ws.Range("I135").Value = 2176.9473
$ws.Range("J135").Value = 1878.6
$ws.Range("K135").Value = 19592.5257
$ws.Range("L135").Value = 16907.4
$ws.Range("M135").Value = -17057.5257
$ws.Range("N135").Value = -21977.4
$ws.Range("H138").Value = 4710.507
$ws.Range("J138").Value = 4911.82
$ws.Range("L138").Value = 14735.46
$ws.Range("N138").Value = -25015.46

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2434.4546
$ws.Range("I110").Value = 2434.4546
$ws.Range("K110").Value = 2434.4546
$ws.Range("M110").Value = -389.4546
$ws.Range("H119").Value = 85227.23
$ws.Range("J119").Value = 85227.23
$ws.Range("L119").Value = 85227.23
$ws.Range("N119").Value = -94903.23
$ws.Range("H135").Value = 64356.082
$ws.Range("J135").Value = 64356.082
$ws.Range("L135").Value = 64356.082
$ws.Range("N135").Value = -74496.08199999999
$ws.Range("H137").Value = 87064.586
$ws.Range("J137").Value = 88745.46000000001
$ws.Range("L137").Value = 88745.46000000001
$ws.Range("N137").Value = -98945.46000000001
$ws.Range("H139").Value = 77143.34
$ws.Range("J139").Value = 77143.34
$ws.Range("L139").Value = 77143.34
$ws.Range("N139").Value = -87423.34

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H88").Value = 39997
$ws.Range("J88").Value = 39997
$ws.Range("L88").Value = 39997
$ws.Range("N88").Value = -40809
$ws.Range("H91").Value = 39997
$ws.Range("J91").Value = 39997
$ws.Range("L91").Value = 39997
$ws.Range("N91").Value = -42805
$ws.Range("H132").Value = 88750.86
$ws.Range("J132").Value = 88750.86
$ws.Range("L132").Value = 88750.86
$ws.Range("N132").Value = -98870.86
$ws.Range("H134").Value = 1842.8064
$ws.Range("I134").Value = 1556.2069
$ws.Range("K134").Value = 4668.620699999999
$ws.Range("M134").Value = -2133.620699999999
$ws.Range("H137").Value = 78403.11
$ws.Range("J137").Value = 78403.11
$ws.Range("L137").Value = 78403.11
$ws.Range("N137").Value = -88603.11

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7169.6787
$ws.Range("I31").Value = 6679.278
$ws.Range("J31").Value = 8052.4
$ws.Range("K31").Value = 6679.278
$ws.Range("L31").Value = 8052.4
$ws.Range("M31").Value = -6384.278
$ws.Range("N31").Value = -8642.4
$ws.Range("H34").Value = 7169.6787
$ws.Range("I34").Value = 6679.278
$ws.Range("J34").Value = 8052.4
$ws.Range("K34").Value = 6679.278
$ws.Range("L34").Value = 8052.4
$ws.Range("M34").Value = -6477.278
$ws.Range("N34").Value = -8456.4
$ws.Range("H114").Value = 74984
$ws.Range("J114").Value = 74984
$ws.Range("L114").Value = 74984
$ws.Range("N114").Value = -83662
$ws.Range("H134").Value = 1193.8518
$ws.Range("I134").Value = 959.7917
$ws.Range("J134").Value = 3066.3333
$ws.Range("K134").Value = 2879.3751
$ws.Range("L134").Value = 9198.999899999999
$ws.Range("M134").Value = -344.3751000000002
$ws.Range("N134").Value = -14268.9999
$ws.Range("H138").Value = 88939.89999999999
$ws.Range("J138").Value = 88939.89999999999
$ws.Range("L138").Value = 88939.89999999999
$ws.Range("N138").Value = -99219.89999999999
$ws.Range("H140").Value = 89775.3
$ws.Range("J140").Value = 89997.42
$ws.Range("L140").Value = 89997.42
$ws.Range("N140").Value = -100357.42
$ws.Range("H141").Value = 229494.22
$ws.Range("J141").Value = 240052.7
$ws.Range("L141").Value = 240052.7
$ws.Range("N141").Value = -250412.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 705.4286
$ws.Range("I29").Value = 729
$ws.Range("J29").Value = 663
$ws.Range("K29").Value = 2187
$ws.Range("L29").Value = 1989
$ws.Range("M29").Value = -1910
$ws.Range("N29").Value = -2543
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H56").Value = 4646.614
$ws.Range("I56").Value = 4646.614
$ws.Range("K56").Value = 4646.614
$ws.Range("M56").Value = -4116.614

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 625
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840
$ws.Range("H132").Value = 5591.1523
$ws.Range("I132").Value = 5239.826
$ws.Range("K132").Value = 15719.478
$ws.Range("M132").Value = -13189.478
$ws.Range("H133").Value = 89860.14999999999
$ws.Range("J133").Value = 89860.14999999999
$ws.Range("L133").Value = 89860.14999999999
$ws.Range("N133").Value = -99980.14999999999
$ws.Range("H135").Value = 87028.34
$ws.Range("J135").Value = 87028.34
$ws.Range("L135").Value = 87028.34
$ws.Range("N135").Value = -97168.34

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 59281.332
$ws.Range("J95").Value = 59281.332
$ws.Range("L95").Value = 59281.332
$ws.Range("N95").Value = -64773.332
$ws.Range("H136").Value = 4673.4736
$ws.Range("I136").Value = 3179.8
$ws.Range("J136").Value = 10274.75
$ws.Range("K136").Value = 9539.400000000001
$ws.Range("L136").Value = 30824.25
$ws.Range("M136").Value = -6989.400000000001
$ws.Range("N136").Value = -35924.25
$ws.Range("H139").Value = 88224
$ws.Range("J139").Value = 88224
$ws.Range("L139").Value = 88224
$ws.Range("N139").Value = -98504
$ws.Range("H141").Value = 83910.31
$ws.Range("J141").Value = 83910.31
$ws.Range("L141").Value = 83910.31
$ws.Range("N141").Value = -94270.31

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H136").Value = 17622.217
$ws.Range("I136").Value = 20650.834
$ws.Range("K136").Value = 61952.50199999999
$ws.Range("M136").Value = -59402.50199999999
$ws.Range("H137").Value = 89425.96000000001
$ws.Range("J137").Value = 89425.96000000001
$ws.Range("L137").Value = 89425.96000000001
$ws.Range("N137").Value = -99625.96000000001
$ws.Range("H139").Value = 87082.46000000001
$ws.Range("I139").Value = 66650
$ws.Range("J139").Value = 87683.414
$ws.Range("K139").Value = 66650
$ws.Range("L139").Value = 87683.414
$ws.Range("M139").Value = -61510
$ws.Range("N139").Value = -97963.414

